$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# Add a new row 5 with a new user entry (mirrors the "ascender usuario" action
# that used to overwrite user/email but no longer does).
$ws.Range("A5").Value = "qwe"
$ws.Range("C5").Value = "zxc"
$ws.Range("B5").Value = "@mail"
$ws.Range("D5").Value = "asd"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 4

$ws.Range("F9").Select()
